$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: extend header row (row 1) formatting from AG1 into AH1:AR1 ---
# Copy the existing header cell format (bold/border/center-top, style index reused)
$ws.Range("AG1").Copy()
$ws.Range("AH1:AR1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Step 2: write the new AH:AR block (columns 34-44) for rows 1-28 ---
$newBlock = @{}
$newBlock[1] = @(50.6, 52, 52.9, 53.2, 57.2, 58.2, 59.2, 61.5, 62, 70, 76.6)
$newBlock[2] = @(0.2, -0.8, -0.1, 0, 0.6, 0.6, 0.6, 0.6, 0.6, 0, 0)
$newBlock[3] = @(0.2, 0.2, 0.1, 0.1, 0.5, 0.5, 0.5, 0.3, 0.2, 0.2, 0)
$newBlock[4] = @(1.1, 0.9, 0.6, 1, 1, 1, 1, 1, 1, 1, -5)
$newBlock[5] = @(0.2, 0.2, 0, 0, 0, 0, 0, 0, 0, 0, 3)
$newBlock[6] = @(0.1, 0, 0.1, 0.1, 0.4, 0.4, 0.4, 0.3, 0.1, 0.1, 1)
$newBlock[7] = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0, -10)
$newBlock[8] = @(3, 3, 3, 0.5, 0.5, 0.5, 0.5, 0.5, 3, 0.1, 0.1)
$newBlock[9] = @(0, 0, 0, 0.4, 4.5, 5.5, 6.6, 13.2, 12.6, 12.6, 12.6)
$newBlock[10] = @(0, 0, 0, 0.5, 4.5, 5.7, 6.1, 6.3, 6.3, 6.3, 6.3)
$newBlock[11] = @(0, 0, 0, -0.8, -0.4, -0.4, -0.2, 0.3, 0, 0, 0)
$newBlock[12] = @(0.2, 0.2, 0, 0, 0, 0, 0, 0, 0, 2, 2)
$newBlock[13] = @(0, 0, 0, 0.13, 0.3, 0.3, 0.2, 0.1, 0, 0, 0)
$newBlock[14] = @(0.1, 0.1, 0, 0, 0, 0, 0, 0.1, 0.1, 0.1, 0.1)
$newBlock[15] = @(-0.1, -0.1, 0, -0.1, -0.1, -0.1, -0.1, -0.1, -0.1, -0.1, -0.1)
$newBlock[16] = @(0.2, 0.2, -0.1, -0.1, -0.1, -0.1, -0.1, 0, 0, 2, 2)
$newBlock[17] = @(-0.02, -0.02, -0.31, 0, 0.3, 0.3, 0.3, 0.1, -0.02, -0.02, -0.02)
$newBlock[18] = @(0, 0, -0.2, 0, 0, 0, 0, 0, 0, 0, 0)
$newBlock[19] = @(0.7, 0.7, 0, 0, 0.5, 0.6, 0.6, 1.1, 0.7, 0.7, 0.7)
$newBlock[20] = @(0.1, 0, 0.43, -0.1, 0.3, 0.3, 0.2, 0.2, 0.1, 0.1, 0.1)
$newBlock[21] = @(0.3, 0.3, 0.2, 0.2, 0.2, 0.2, 0.2, 0.2, 0.2, 0.2, 0.2)
$newBlock[22] = @(0.2, 0, 0.8, 0.1, 0.6, 0.7, 0.8, 0.6, 0.2, 0.2, 0.2)
$newBlock[23] = @(0.1, -0.1, 1, -1, 0, 0, 0, 0.4, 0.2, 0.2, 0.2)
$newBlock[24] = @(-0.5, -0.5, 0.1, 0, 0, 0, 0, 0, -0.5, -0.5, -0.5)
$newBlock[25] = @(0.1, 0.1, 0.4, 0.5, 0.4, 0.2, -0.1, 0.2, 0.1, 0.1, 0.1)
$newBlock[26] = @(-0.1, -0.1, 0.1, 0, 0.4, 0.3, 0.3, 0.1, 0.1, 0.1, 0.1)
$newBlock[27] = @(0.1, 0.1, 0.8, 0.9, 0.1, 0.5, 0.1, 0.3, 0.1, 0.1, 0.1)
$newBlock[28] = @(-0.1, -0.1, 1.1, 0, 0, -0.3, -0.5, -0.8, -0.1, -0.1, -0.1)

foreach ($r in $newBlock.Keys) {
    $rowVals = $newBlock[$r]
    for ($i = 0; $i -lt $rowVals.Length; $i++) {
        $ws.Cells.Item($r, 34 + $i).Value = $rowVals[$i]
    }
}

# --- Step 3: apply the scattered single-cell edits within the existing B:AG block ---
$singleEdits = @{
    "AE2" = 0.4
    "AE3" = 0.2
    "AE4" = 1.1
    "N6" = 0.3
    "N15" = -0.1
    "N17" = 0
    "N19" = -0.3
    "O19" = 0.3
    "N20" = -0.1
    "O20" = 0.1
    "N22" = 0.1
    "O22" = 0.6
    "N23" = -1
    "N28" = 0
}
foreach ($addr in $singleEdits.Keys) {
    $ws.Range($addr).Value = $singleEdits[$addr]
}

Write-Host "edit complete"
